# Apply the "add marking PNI data to enhancement folder" edit:
#  - rename the existing "Sheet1" to "old"
#  - add a new worksheet "PBT" after it, with brood-year / PBT marking status data
#  - adjust sheet selections / active sheet to match the authored workbook

$wb = $excel.ActiveWorkbook

# --- rename the original sheet ---
$old = $wb.Worksheets.Item(1)
$old.Name = "old"

# clear any stale selection on the old sheet and select A1:G1 (header row)
$old.Activate() | Out-Null
$old.Range("A1:G1").Select() | Out-Null

# --- add the new "PBT" sheet right after "old" ---
$pbt = $wb.Worksheets.Add($null, $old)
$pbt.Name = "PBT"

# header row
$pbt.Range("A1").Value = "BROOD_YEAR"
$pbt.Range("B1").Value = "PBT"
$pbt.Range("C1").Value = "comments"

# data rows
$pbt.Range("A2").Value = 2017
$pbt.Range("B2").Value = "Incomplete"

$pbt.Range("A3").Value = 2018
$pbt.Range("B3").Value = "Partial*"
$pbt.Range("C3").Value = "*Should be full soon"

$pbt.Range("A4").Value = 2019
$pbt.Range("B4").Value = "Full"

$pbt.Range("A5").Value = 2020
$pbt.Range("B5").Value = "Full"

$pbt.Range("A6").Value = 2021
$pbt.Range("B6").Value = "Full"

$pbt.Range("A7").Value = 2022
$pbt.Range("B7").Value = "Full"

# column A sizing (bestFit on the original), closest achievable width
$pbt.Columns.Item(1).ColumnWidth = 12

# print setup to match the authored sheet
$pbt.PageSetup.Orientation = 1

# select / activate like the saved workbook (PBT tab active, cell I11 selected)
$pbt.Range("I11").Select() | Out-Null
$pbt.Activate() | Out-Null

Write-Output "done"
